$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.093.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.264.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.41%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '397.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.02%  '

$ws.Range("E7").Value = '  +4.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.623'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("E10").Value = '  +0.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0957'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.24%  '

$ws.Range("E12").Value = '  +1.96%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.778.88'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.28'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.256.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '56.920.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.03%  '

$ws.Range("E20").Value = '  -0.65%  '

$ws.Range("E21").Value = '  +6.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.94'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '293.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.30%  '

$ws.Range("E27").Value = '  -3.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.16%  '

$ws.Range("E30").Value = '  -2.31%  '

$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.112'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.69%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '39.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0486'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.67%  '

$ws.Range("E36").Value = '  +0.97%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.09%  '

$ws.Range("E39").Value = '  -1.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '136.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.48%  '

$ws.Range("E42").Value = '  +1.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.284'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.03%  '

$ws.Range("E44").Value = '  -2.41%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.82%  '

$ws.Range("E46").Value = '  -1.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.59%  '

$ws.Range("E48").Value = '  +4.69%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.152.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.99'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.64%  '

$ws.Range("E51").Value = '  -5.80%  '
